# MENU_MOCK.xlsx edit: remove the OptionBluntMultiplier row (row 11) from the
# "Menu Options" sheet. Deleting the entire row shifts every row below it up
# by one, which is exactly what the target diff shows (OptionFireMultiplier
# moves from row 12 -> 11, ... OptionResetStats moves from row 57 -> 56, and
# the very last row disappears). The sheet's used range shrinks from
# A1:G57 to A1:G56 automatically as part of the row deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 holds: CategoryDamageTypeMultipliers | 30 | OptionBluntMultiplier |
# float | 0.5f | "DOT damage multiplier for blunt attacks. ..." | (blank)
$ws.Rows("11").Delete()
